# Generate Report for Handback
# Update the handoff/handback timestamp cells across the report sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value = "2016-09-01 13:14:43"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-01 13:14:38"
$wsZhCn.Range("K2").Value = "2016-09-01 13:14:59"

# de-de sheet: Correspond Handoff Datetime (mirrors the Overview value, since
# both cells originally shared the same underlying string) and
# Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-09-01 13:14:43"
$wsDeDe.Range("K2").Value = "2016-09-01 13:15:28"
